# Update the "想去人数" (want-to-go count) column F on each sheet to the
# freshly scraped values (output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 7418
    4  = 3539
    6  = 3864
    8  = 89
    10 = 107
    11 = 161
    14 = 147
    19 = 4170
    22 = 1031
    24 = 1878
    26 = 98
    27 = 3069
    28 = 2291
    30 = 85
    32 = 47
    33 = 115
    34 = 43
    36 = 4371
    37 = 494
    38 = 326
    41 = 827
    42 = 227
    44 = 1653
    47 = 616
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    4  = 442
    16 = 600
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "本地生活" ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 168
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 168
    5  = 7418
    6  = 3539
    7  = 3864
    9  = 89
    11 = 107
    13 = 161
    16 = 147
    21 = 4170
    27 = 1878
    29 = 98
    30 = 3069
    31 = 2291
    33 = 85
    34 = 116
    37 = 4371
    39 = 494
    40 = 326
    42 = 827
    43 = 227
    45 = 1653
    48 = 616
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
